$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-19.14024064054204, 1.96079183533252, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, 1.321359705725766, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, 1.993050230719507, 1.624259248992513, -19.14024064054204, 3.401035856060555, -19.14024064054204, 1.68383630854932, -19.14024064054204, 0.7859525441903138, -19.14024064054204),
    @(-19.14024064054204, 1.665645686933433, -19.14024064054204, -19.14024064054204, -19.14024064054204, 2.750690437242351, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204),
    @(2.474782304374543, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, -19.14024064054204, 1.784546497169117, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204),
    @(3.852112213379589, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, 1.801177412603078, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, -19.14024064054204, 2.907726290294515, -19.14024064054204, 2.765108534252993, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, -19.14024064054204, 2.531366506485722, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, 1.56719685425178, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, 1.557542279008511, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, 4.32192572290233),
    @(-19.14024064054204, -19.14024064054204, 1.768513318309847, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, 1.841893009386286, -19.14024064054204),
    @(-19.14024064054204, 2.143035965714628, 1.847371038906726, -19.14024064054204, -19.14024064054204, -19.14024064054204, 2.175599373996466, 2.021986555204631, 2.606182239788768, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, -19.14024064054204, 1.802065524233347, 2.066895849075365, 2.495376968372725, -19.14024064054204),
    @(-19.14024064054204, -19.14024064054204, 2.048665078240776, -19.14024064054204, -19.14024064054204, -19.14024064054204, 1.522355594247933, 1.809720684095673, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, 1.039086561090437, 1.50127226166072, -19.14024064054204, 3.238224890624981, -19.14024064054204, 1.671645926459797, 1.175622033362656, -19.14024064054204, -19.14024064054204),
    @(-19.14024064054204, 1.311534919951295, -19.14024064054204, 1.708997217579143, -19.14024064054204, 2.694119685286902, 1.445654392626248, -19.14024064054204, -19.14024064054204, -19.14024064054204)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
